# Weekly price update: insert a new data row for the most recent week
# (2021-09-10, serial 44449) above the existing "Feria Lagunitas de Puerto
# Montt - Arveja Verde / Perfection" block, pushing the previous rows
# (formerly 59-68) down to 60-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 59; this shifts rows 59:68 down to 60:69
# and extends the used range to A1:R69 automatically.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with this week's record.
$ws.Range("A59").Value = 4
$ws.Range("B59").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C59").Value = "Los Lagos"
$ws.Range("D59").Value = 44449
$ws.Range("E59").Value = 10
$ws.Range("F59").Value = 100112022
$ws.Range("G59").Value = "Arveja Verde"
$ws.Range("H59").Value = "Perfection"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 70
$ws.Range("K59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("M59").Value = 40000
$ws.Range("N59").Value = "`$/malla 25 kilos"
$ws.Range("O59").Value = "Provincia de Huasco"
$ws.Range("P59").Value = 1600
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"
